# Generate Report for Archive
# Update the localization status report: mark the two files that are
# currently mid-flight (9121ae3b... and a30f3e4a...) as "In Translation"
# instead of "Ready for handoff" across the Overview rollup sheet and the
# per-locale detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de) for rows 7-8 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus
$overview.Range("B8").Value = $newStatus
$overview.Range("C8").Value = $newStatus

# --- zh-cn detail sheet: Status column (C) for rows 7-8 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = $newStatus
$zhcn.Range("C8").Value = $newStatus

# --- de-de detail sheet: Status column (C) for rows 7-8 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = $newStatus
$dede.Range("C8").Value = $newStatus
